$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2,2,2,10,0,2,10,0,7,7,2,2,2,2,2,2),
    @(13,5,5,12,1,5,12,1,5,9,5,5,5,5,5,5),
    @(1,7,1,7,2,1,7,2,9,2,1,10,8,10,1,10),
    @(5,1,13,2,3,13,2,3,8,5,13,1,10,8,10,8),
    @(7,3,7,9,4,7,9,4,6,8,7,8,1,1,7,1),
    @(3,0,10,5,5,10,5,5,2,0,10,6,7,6,13,6),
    @(10,3,3,1,6,3,1,6,0,4,3,7,13,7,3,7),
    @(12,2,12,8,7,12,8,7,4,1,12,9,6,0,12,9),
    @(9,9,9,3,8,9,3,8,3,6,9,0,3,9,8,0),
    @(8,8,8,13,9,8,13,9,1,3,8,11,0,11,9,11),
    @(6,4,6,0,10,6,0,10,0,0,0,14,9,12,6,14),
    @(0,0,0,14,11,0,14,11,1,1,6,4,12,13,0,3),
    @(14,6,14,6,12,14,6,12,1,2,14,3,11,14,14,12),
    @(11,1,11,4,13,11,4,13,5,3,11,12,14,3,11,4),
    @(4,4,4,11,14,4,11,14,3,0,13,13,4,4,4,13)
)

for ($r = 0; $r -lt $values.Count; $r++) {
    $rowVals = $values[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item(2 + $r, 2 + $c).Value = $rowVals[$c]
    }
}
